$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New tracklist data (Elle Varner "ellevarner2" set) replacing the old
# "common9" (Common) tracklist, for both Sheet1 and Sheet3 (the two
# sheets that hold the raw web-query results behind the formatted
# report on Sheet2).
# ---------------------------------------------------------------------

$titles = @(
    "Coffee on the Roof",
    "Pour Me",
    "1 to 10",
    "Wishing Well",
    "Number One Song",
    "Loving U Blind",
    "Kinda Love",
    "Casanova",
    "Be Encouraged"
)

$composers = @(
    "Elle Varner / Jimmy Varner",
    "Olubowale Akintimehin / Elle Varner / Jimmy Varner",
    "Elle Varner / Jordan Ware",
    "Marlanna Evans / Elle Varner / Jimmy Varner",
    "Elle Varner / Jimmy Varner",
    "Elle Varner / Jimmy Varner",
    "Stacy Barthe / Los Hendrix / L3gion / Elle Varner / Jordan Ware",
    "Elle Varner",
    "Nascent / Elle Varner / Jimmy Varner"
)

$performers = @(
    "Elle Varner",
    "Elle Varner feat. Wale",
    "Elle Varner",
    "Elle Varner feat. Rapsody",
    "Elle Varner",
    "Elle Varner",
    "Elle Varner",
    "Elle Varner",
    "Elle Varner"
)

$times = @(
    0.18611111111111112,
    0.17083333333333331,
    0.15694444444444444,
    0.18055555555555555,
    0.21249999999999999,
    0.14444444444444446,
    0.13472222222222222,
    0.17430555555555557,
    0.13125000000000001
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row: merge "Title"/"Composer" into a single "Title/Composer"
    # header living in column B; drop the old " No." (A1) / "Composer" (C1)
    # header text entirely.
    $ws.Range("A1").ClearContents()
    $ws.Range("B1").Value = "Title/Composer"
    $ws.Range("C1").ClearContents()
    # D1 (Performer) / E1 (Time) stay as-is.

    # Rows 2-10: the 9 new tracks.
    for ($i = 0; $i -lt $titles.Length; $i++) {
        $r = $i + 2
        $ws.Range("A$r").Value = ($i + 1)
        $ws.Range("B$r").Value = $titles[$i]
        $ws.Range("C$r").Value = $composers[$i]
        $ws.Range("D$r").Value = $performers[$i]
        $ws.Range("E$r").Value = $times[$i]
    }

    # Rows 11-12 no longer hold data (old tracks 10 & 11 removed) - only
    # the time column's formatting placeholder cell remains.
    $ws.Range("A11:E12").ClearContents()

    # Column widths shrink to fit the now-shorter content.
    $ws.Columns.Item(1).ColumnWidth = 1.8571428571428572
    $ws.Columns.Item(2).ColumnWidth = 17.285714285714285
    $ws.Columns.Item(3).ColumnWidth = 58.57142857142857
    $ws.Columns.Item(4).ColumnWidth = 23
}

# ---------------------------------------------------------------------
# Rename the defined names (and, implicitly, the underlying query/page
# slug) from "common9" to "ellevarner2".
# ---------------------------------------------------------------------
$wb.Names.Item("Sheet1!common9").Name = "ellevarner2"
$wb.Names.Item("Sheet3!common9").Name = "ellevarner2"

# ---------------------------------------------------------------------
# Sheet2 (the formatted report) recalculates automatically from the new
# Sheet1 data; just restore the selection/active-cell state recorded in
# the workbook.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K13").Select()
